$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (timestamp values already in Excel serial date form)
$newData = @(
    @{ Row = 46; A = 5493; B = 45741.45833333334 },
    @{ Row = 47; A = 5468; B = 45741.46875 },
    @{ Row = 48; A = 5553; B = 45741.47916666666 }
)

foreach ($entry in $newData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    # Match the existing timestamp column formatting used by the rows above
    $ws.Cells.Item($r, 2).NumberFormat = $ws.Cells.Item($r - 1, 2).NumberFormat
}

